$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AeroToolKitFunctionList")
$win = $excel.ActiveWindow
Write-Output $win.Zoom
